$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the existing task description (in place, keeps shared-string index)
[void]$ws.Cells.Replace("my Marting", "my Martin")

# Correct the dates for the two existing "gallery" rows (2008 -> 2018)
$ws.Range("B16").Value = 43361
$ws.Range("B17").Value = 43362

# Insert two new rows before the old total row (currently row 19) so the
# total moves down to row 21, with a blank separator row at row 20
[void]$ws.Rows("18:19").Insert()

$ws.Range("A18").Value = "Added 2008 Images and prices"
$ws.Range("B18").Value = 43363
$ws.Range("C18").Value = 2

$ws.Range("A19").Value = "Added remaining graphics"
$ws.Range("B19").Value = 43364
$ws.Range("C19").Value = 1.5

# Update the active selection to match the new state
[void]$ws.Range("B20").Select()
